$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = 'ECs'
$ws.Cells.Item(2,2).Value = 'Il19'
$ws.Cells.Item(2,3).Value = 'Il20rb'
$ws.Cells.Item(2,4).Value = 'ECs'
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.5
$ws.Cells.Item(2,7).Value = 0.0621735
$ws.Cells.Item(2,8).Value = 0.124347
$ws.Cells.Item(2,9).Value = 0.400116482236459
$ws.Cells.Item(2,10).Value = 0.3077956989247312
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 5.3589275
$ws.Cells.Item(2,14).Value = 10.717855
$ws.Cells.Item(2,15).Value = 0.3905787093313861
$ws.Cells.Item(2,16).Value = 0.3109096775151767
$ws.Cells.Item(2,17).Value = 0.33318327892125
$ws.Cells.Item(2,18).Value = 1.332733115685
$ws.Cells.Item(2,19).Value = 0.1562769792141306
$ws.Cells.Item(2,20).Value = 0.09569666149324658

# Row 3
$ws.Cells.Item(3,1).Value = 'ECs'
$ws.Cells.Item(3,2).Value = 'Il19'
$ws.Cells.Item(3,3).Value = 'Il20rb'
$ws.Cells.Item(3,4).Value = 'FAPs'
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.5
$ws.Cells.Item(3,7).Value = 0.0621735
$ws.Cells.Item(3,8).Value = 0.124347
$ws.Cells.Item(3,9).Value = 0.400116482236459
$ws.Cells.Item(3,10).Value = 0.3077956989247312
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 2.206973333333333
$ws.Cells.Item(3,14).Value = 6.62092
$ws.Cells.Item(3,15).Value = 0.1608524832743344
$ws.Cells.Item(3,16).Value = 0.192063440124333
$ws.Cells.Item(3,17).Value = 0.13721525654
$ws.Cells.Item(3,18).Value = 0.82329153924
$ws.Cells.Item(3,19).Value = 0.06435972976672551
$ws.Cells.Item(3,20).Value = 0.05911630079095735

# Row 4
$ws.Cells.Item(4,1).Value = 'ECs'
$ws.Cells.Item(4,2).Value = 'Il19'
$ws.Cells.Item(4,3).Value = 'Il20rb'
$ws.Cells.Item(4,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.5
$ws.Cells.Item(4,7).Value = 0.0621735
$ws.Cells.Item(4,8).Value = 0.124347
$ws.Cells.Item(4,9).Value = 0.400116482236459
$ws.Cells.Item(4,10).Value = 0.3077956989247312
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.315162333333333
$ws.Cells.Item(4,14).Value = 3.945487
$ws.Cells.Item(4,15).Value = 0.09585395710514605
$ws.Cells.Item(4,16).Value = 0.1144529470505359
$ws.Cells.Item(4,17).Value = 0.0817682453315
$ws.Cells.Item(4,18).Value = 0.490609471989
$ws.Cells.Item(4,19).Value = 0.03835274812535547
$ws.Cells.Item(4,20).Value = 0.03522812483141496

# Row 5
$ws.Cells.Item(5,1).Value = 'ECs'
$ws.Cells.Item(5,2).Value = 'Il19'
$ws.Cells.Item(5,3).Value = 'Il20rb'
$ws.Cells.Item(5,4).Value = 'MuSCs'
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.5
$ws.Cells.Item(5,7).Value = 0.0621735
$ws.Cells.Item(5,8).Value = 0.124347
$ws.Cells.Item(5,9).Value = 0.400116482236459
$ws.Cells.Item(5,10).Value = 0.3077956989247312
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.3299455
$ws.Cells.Item(5,14).Value = 2.659891
$ws.Cells.Item(5,15).Value = 0.09693140966566258
$ws.Cells.Item(5,16).Value = 0.0771596418346321
$ws.Cells.Item(5,17).Value = 0.08268736654425
$ws.Cells.Item(5,18).Value = 0.330749466177
$ws.Cells.Item(5,19).Value = 0.038783854653646
$ws.Cells.Item(5,20).Value = 0.02374940588727251

# Row 6
$ws.Cells.Item(6,1).Value = 'ECs'
$ws.Cells.Item(6,2).Value = 'Il19'
$ws.Cells.Item(6,3).Value = 'Il20rb'
$ws.Cells.Item(6,4).Value = 'Neutrophils'
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.5
$ws.Cells.Item(6,7).Value = 0.0621735
$ws.Cells.Item(6,8).Value = 0.124347
$ws.Cells.Item(6,9).Value = 0.400116482236459
$ws.Cells.Item(6,10).Value = 0.3077956989247312
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 2.739463
$ws.Cells.Item(6,14).Value = 8.218389
$ws.Cells.Item(6,15).Value = 0.1996623247470855
$ws.Cells.Item(6,16).Value = 0.2384037359792865
$ws.Cells.Item(6,17).Value = 0.1703220028305
$ws.Cells.Item(6,18).Value = 1.021932016983
$ws.Cells.Item(6,19).Value = 0.07988818701295734
$ws.Cells.Item(6,20).Value = 0.07337964454201155

# Row 7
$ws.Cells.Item(7,1).Value = 'ECs'
$ws.Cells.Item(7,2).Value = 'Il19'
$ws.Cells.Item(7,3).Value = 'Il20rb'
$ws.Cells.Item(7,4).Value = 'Resolving-Mac'
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.5
$ws.Cells.Item(7,7).Value = 0.0621735
$ws.Cells.Item(7,8).Value = 0.124347
$ws.Cells.Item(7,9).Value = 0.400116482236459
$ws.Cells.Item(7,10).Value = 0.3077956989247312
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.7700086666666666
$ws.Cells.Item(7,14).Value = 2.310026
$ws.Cells.Item(7,15).Value = 0.05612111587638537
$ws.Cells.Item(7,16).Value = 0.06701055749603567
$ws.Cells.Item(7,17).Value = 0.047874133837
$ws.Cells.Item(7,18).Value = 0.287244803022
$ws.Cells.Item(7,19).Value = 0.022454983463644
$ws.Cells.Item(7,20).Value = 0.02062556137982818

# Row 8
$ws.Cells.Item(8,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(8,2).Value = 'Il19'
$ws.Cells.Item(8,3).Value = 'Il20rb'
$ws.Cells.Item(8,4).Value = 'ECs'
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = 0.3333333333333333
$ws.Cells.Item(8,7).Value = 0.09321499999999999
$ws.Cells.Item(8,8).Value = 0.279645
$ws.Cells.Item(8,9).Value = 0.5998835177635411
$ws.Cells.Item(8,10).Value = 0.6922043010752689
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 5.3589275
$ws.Cells.Item(8,14).Value = 10.717855
$ws.Cells.Item(8,15).Value = 0.3905787093313861
$ws.Cells.Item(8,16).Value = 0.3109096775151767
$ws.Cells.Item(8,17).Value = 0.4995324269124999
$ws.Cells.Item(8,18).Value = 2.997194561475
$ws.Cells.Item(8,19).Value = 0.2343017301172555
$ws.Cells.Item(8,20).Value = 0.2152130160219301

# Row 9
$ws.Cells.Item(9,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(9,2).Value = 'Il19'
$ws.Cells.Item(9,3).Value = 'Il20rb'
$ws.Cells.Item(9,4).Value = 'FAPs'
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 0.3333333333333333
$ws.Cells.Item(9,7).Value = 0.09321499999999999
$ws.Cells.Item(9,8).Value = 0.279645
$ws.Cells.Item(9,9).Value = 0.5998835177635411
$ws.Cells.Item(9,10).Value = 0.6922043010752689
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 2.206973333333333
$ws.Cells.Item(9,14).Value = 6.62092
$ws.Cells.Item(9,15).Value = 0.1608524832743344
$ws.Cells.Item(9,16).Value = 0.192063440124333
$ws.Cells.Item(9,17).Value = 0.2057230192666667
$ws.Cells.Item(9,18).Value = 1.8515071734
$ws.Cells.Item(9,19).Value = 0.09649275350760884
$ws.Cells.Item(9,20).Value = 0.1329471393333757

# Row 10
$ws.Cells.Item(10,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(10,2).Value = 'Il19'
$ws.Cells.Item(10,3).Value = 'Il20rb'
$ws.Cells.Item(10,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0.3333333333333333
$ws.Cells.Item(10,7).Value = 0.09321499999999999
$ws.Cells.Item(10,8).Value = 0.279645
$ws.Cells.Item(10,9).Value = 0.5998835177635411
$ws.Cells.Item(10,10).Value = 0.6922043010752689
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.315162333333333
$ws.Cells.Item(10,14).Value = 3.945487
$ws.Cells.Item(10,15).Value = 0.09585395710514605
$ws.Cells.Item(10,16).Value = 0.1144529470505359
$ws.Cells.Item(10,17).Value = 0.1225928569016667
$ws.Cells.Item(10,18).Value = 1.103335712115
$ws.Cells.Item(10,19).Value = 0.05750120897979059
$ws.Cells.Item(10,20).Value = 0.07922482221912099

# Row 11
$ws.Cells.Item(11,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(11,2).Value = 'Il19'
$ws.Cells.Item(11,3).Value = 'Il20rb'
$ws.Cells.Item(11,4).Value = 'MuSCs'
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 0.3333333333333333
$ws.Cells.Item(11,7).Value = 0.09321499999999999
$ws.Cells.Item(11,8).Value = 0.279645
$ws.Cells.Item(11,9).Value = 0.5998835177635411
$ws.Cells.Item(11,10).Value = 0.6922043010752689
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 1.3299455
$ws.Cells.Item(11,14).Value = 2.659891
$ws.Cells.Item(11,15).Value = 0.09693140966566258
$ws.Cells.Item(11,16).Value = 0.0771596418346321
$ws.Cells.Item(11,17).Value = 0.1239708697825
$ws.Cells.Item(11,18).Value = 0.743825218695
$ws.Cells.Item(11,19).Value = 0.05814755501201658
$ws.Cells.Item(11,20).Value = 0.05341023594735959

# Row 12
$ws.Cells.Item(12,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(12,2).Value = 'Il19'
$ws.Cells.Item(12,3).Value = 'Il20rb'
$ws.Cells.Item(12,4).Value = 'Neutrophils'
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(12,6).Value = 0.3333333333333333
$ws.Cells.Item(12,7).Value = 0.09321499999999999
$ws.Cells.Item(12,8).Value = 0.279645
$ws.Cells.Item(12,9).Value = 0.5998835177635411
$ws.Cells.Item(12,10).Value = 0.6922043010752689
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 2.739463
$ws.Cells.Item(12,14).Value = 8.218389
$ws.Cells.Item(12,15).Value = 0.1996623247470855
$ws.Cells.Item(12,16).Value = 0.2384037359792865
$ws.Cells.Item(12,17).Value = 0.255359043545
$ws.Cells.Item(12,18).Value = 2.298231391905
$ws.Cells.Item(12,19).Value = 0.1197741377341282
$ws.Cells.Item(12,20).Value = 0.1650240914372749

# Row 13
$ws.Cells.Item(13,1).Value = 'Inflammatory-Mac'
$ws.Cells.Item(13,2).Value = 'Il19'
$ws.Cells.Item(13,3).Value = 'Il20rb'
$ws.Cells.Item(13,4).Value = 'Resolving-Mac'
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 0.3333333333333333
$ws.Cells.Item(13,7).Value = 0.09321499999999999
$ws.Cells.Item(13,8).Value = 0.279645
$ws.Cells.Item(13,9).Value = 0.5998835177635411
$ws.Cells.Item(13,10).Value = 0.6922043010752689
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.7700086666666666
$ws.Cells.Item(13,14).Value = 2.310026
$ws.Cells.Item(13,15).Value = 0.05612111587638537
$ws.Cells.Item(13,16).Value = 0.06701055749603567
$ws.Cells.Item(13,17).Value = 0.07177635786333332
$ws.Cells.Item(13,18).Value = 0.6459872207699998
$ws.Cells.Item(13,19).Value = 0.03366613241274137
$ws.Cells.Item(13,20).Value = 0.04638499611620749
